$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: player names and "Point:" headers
$ws.Range("A1").Value = "Kim"
$ws.Range("B1").Value = "Point:"
$ws.Range("C1").Value = "Mads"
$ws.Range("D1").Value = "Point:"
$ws.Range("E1").Value = "Soren"
$ws.Range("F1").Value = "Point:"
$ws.Range("G1").Value = "Emil"
$ws.Range("H1").Value = "Point:"

# Row 2: team names and point formulas
$ws.Range("A2").Value = "Chelsea"
$ws.Range("C2").Value = "Arsenal"
$ws.Range("E2").Value = "Manchester Utd"
$ws.Range("G2").Value = "Leicester City"
$ws.Range("B2").Formula = "=0+5+10+5+10+10+5+5+5+5+10+5+20+30+10+10+5+10"
$ws.Range("D2").Formula = "=0+15+30+30+5+5+25+20+10+10+5+15+20+10+10+40+15"
$ws.Range("F2").Formula = "=0+5+10+5+30+30+15+20+10+5+10+5+5+5+5+20+10+10+25+30+10+25+10"
$ws.Range("H2").Formula = "=0+20+10+10+5+5+30+5+40+5+10+20+0+5+15+5+10+30+10+10+5+5+0+10+10"

# Row 3: team names and point formulas
$ws.Range("A3").Value = "RB Leipzig"
$ws.Range("C3").Value = "Eint Frankfurt"
$ws.Range("E3").Value = "Leverkusen"
$ws.Range("G3").Value = "Tottenham"
$ws.Range("B3").Formula = "=0+10+10+20+5+5+10+30+30+10+5+15+10+5+10+10+15+5"
$ws.Range("D3").Formula = "=0+40+5+5+5+5+5+10+15+10+0+20+5+15+15+10+10+10+5+10+15+0+0+30+5+5"
$ws.Range("F3").Formula = "=0+5+20+25+5+15+5+40+10+10+5+10+5+10+5+20"
$ws.Range("H3").Formula = "=0+20+40+30+10+40+5+5+5+30+10+15+10+20+10+5+5"

# Row 4: team names and point formulas
$ws.Range("A4").Value = "Barcelona"
$ws.Range("C4").Value = "Hoffenheim"
$ws.Range("E4").Value = "Sevilla"
$ws.Range("G4").Value = "Dortmund"
$ws.Range("B4").Formula = "=0+5+5+5+15+10+10+5+5+10+5+10+5+5+10+10+5+15"
$ws.Range("D4").Formula = "=0+5+20+15+5+15+25+15+10+5+10+20+15+5+20+10+40+5+10+30+25"
$ws.Range("F4").Formula = "=0+5+10+10+5+5+10+10+10+5+5+5+5+5+10+20+10+5+5+5+5"
$ws.Range("H4").Formula = "=0+10+10+20+10+5+10+40+5+5+40+15+10"

# Row 5: team names and point formulas
$ws.Range("A5").Value = "Bologna"
$ws.Range("C5").Value = "Valencia"
$ws.Range("E5").Value = "Juventus"
$ws.Range("G5").Value = "Real Sociedad"
$ws.Range("B5").Formula = "=0+10+30+5+15+5+30+20+10+10+20+30+10+15+10+5+20+5+10+10+20+10+10+5+5+10+15"
$ws.Range("D5").Formula = "=0+5+10+30+5+5+30+5+20+5+10+5+10+20+10+10+10+10+40+5+5+5+10+15+5+5+20+5"
$ws.Range("F5").Formula = "=0+5+10+10+10+0+5+10+10+20+5+5+10+10+0+0+10+10+10+5+15"
$ws.Range("H5").Formula = "=0+30+10+5+5+5+10+10+15+25+15+5+5+10+25+20+10+5+20+5+10+10"

# Row 6: team names and point formulas
$ws.Range("A6").Value = "AGF"
$ws.Range("C6").Value = "Milan"
$ws.Range("E6").Value = "Torino"
$ws.Range("G6").Value = "Atalanta"
$ws.Range("B6").Formula = "=0+10+5+10+15+30+15+0+0+50+10+0+0+10+5+10+10+10+5+5+5+15+10+0+0+5+10+10+5"
$ws.Range("D6").Formula = "=0+10+5+10+15+5+10+10+10+5+5+10+10"
$ws.Range("F6").Formula = "=0+20+10+5+5+10+20+10+10+5+5+10+5+15+10+10+10+5+10+10+5+10+10+20"
$ws.Range("H6").Formula = "=0+10+10+5+20+5+5+20+5+5+5+10+10+10+10+5+15+10+10+10+5+30+10"

# Row 7: team names and point formulas
$ws.Range("A7").Value = "Odense"
$ws.Range("C7").Value = "FC Copenhagen"
$ws.Range("E7").Value = "Midtjylland"
$ws.Range("G7").Value = "Brøndby"
$ws.Range("B7").Formula = "=0+5+30+15+5+10+15+5+10+5+20+10+30+10+5+10+15+5+5+5+5+5+10"
$ws.Range("D7").Formula = "=0+5+5+20+5+5+20+5+10+15+10+10+10"
$ws.Range("F7").Formula = "=0+20+30+5+40+20+5+15+10+20+10+10+10"
$ws.Range("H7").Formula = "=0+10+5+5+30+10+10+5+10+5+10+20+10+20+10+15+30+10+10+10"

# Row 8: totals (A/C/E/G unchanged text already, ensure set; B/D/F/H sums)
$ws.Range("A8").Value = "Total:"
$ws.Range("C8").Value = "Total:"
$ws.Range("E8").Value = "Total:"
$ws.Range("G8").Value = "Total:"
$ws.Range("B8").Formula = "=SUM(B2:B7)"
$ws.Range("D8").Formula = "=SUM(D2:D7)"
$ws.Range("F8").Formula = "=SUM(F2:F7)"
$ws.Range("H8").Formula = "=SUM(H2:H7)"
